$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "week 3 assignment packages" - the Data-Wrangling / Data-Visualization topics
# (and their "due" labels) were in the wrong order relative to the schedule;
# swap week 4 <-> week 5 content so Data-Visualization comes before Data-Wrangling.

# Swap the "Data-Wrangling" (C7) and "Data-Visualization" (C8) topic cells
$ws.Range("C7").Value = "[Data-Visualization](https://crumplab.github.io/psyc7709/Schedule.html#4_data-visualization)"
$ws.Range("C8").Value = "[Data-Wrangling](https://crumplab.github.io/psyc7709/Schedule.html#5_data-wrangling)"

# Swap the corresponding "due" assignment labels in column D
$ws.Range("D8").Value = "week 4 due (data-vis)"
$ws.Range("D9").Value = "week 5 due (data-wrangling)"

# Update the sheet view: drop the old scroll anchor and move the selection to D10
$ws.Range("D10").Select()

# Reposition/resize the workbook window to match the saved view state
$excel.Left = 80
$excel.Top = 460
$excel.Width = 23220
$excel.Height = 15540
